# Update MRE (marine renewable energy) turbine prices on the "mre_costs" sheet.
# - RM1 total cost (F9): 35561 -> 4360000
# - RM2 rated power (B10): 89.51 -> 90 kW
# - RM2 total cost (F10): 3189 -> 580000
# - RM4 total cost (F11): 102500 -> 15000000

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mre_costs")

$ws.Range("F9").Value = 4360000
$ws.Range("B10").Value = 90
$ws.Range("F10").Value = 580000
$ws.Range("F11").Value = 15000000

# Move the active selection to match the post-edit workbook state.
$ws.Activate()
$ws.Range("B11").Select()
